# Apply updated cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'36.932.66"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.27%  '

$ws.Range('D3').Value = "'2.044.62"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.82%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = "'251.36"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.02%  '

$ws.Range('D6').Value = "'0.669"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.67%  '

$ws.Range('D7').Value = "'58.60"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.80%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('D9').Value = "'61.26"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.79%  '

$ws.Range('E10').Value = '  +1.35%  '

$ws.Range('D11').Value = "'0.0786"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.26%  '

$ws.Range('E12').Value = '  +1.96%  '

$ws.Range('D13').Value = "'16.29"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.79%  '

$ws.Range('D14').Value = "'2.344.08"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.79%  '

$ws.Range('D15').Value = "'0.807"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.38%  '

$ws.Range('D16').Value = "'5.59"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.58%  '

$ws.Range('D17').Value = "'2.045.15"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.77%  '

$ws.Range('D18').Value = "'36.911.32"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.20%  '

$ws.Range('D19').Value = "'16.88"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +15.62%  '

$ws.Range('D20').Value = "'74.89"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.12%  '

$ws.Range('D21').Value = "'0.0₃0905"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.01%  '

$ws.Range('D22').Value = "'5.42"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.94%  '

$ws.Range('D23').Value = "'236.97"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.65%  '

$ws.Range('E24').Value = '  +0.00%  '

$ws.Range('E25').Value = '  -3.89%  '

$ws.Range('D26').Value = "'2.30"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.59%  '

$ws.Range('D27').Value = "'168.96"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.63%  '

$ws.Range('D28').Value = "'9.26"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.17%  '

$ws.Range('D29').Value = "'20.22"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.20%  '

$ws.Range('E30').Value = '  +1.22%  '

$ws.Range('D31').Value = "'1.15"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.02%  '

$ws.Range('D32').Value = "'4.74"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.81%  '

$ws.Range('D33').Value = "'0.0620"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.53%  '

$ws.Range('D34').Value = "'4.47"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.79%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.0881"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.91%  '

$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = "'1.00"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.06%  '

$ws.Range('D37').Value = "'2.22"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.06%  '

$ws.Range('E38').Value = '  -4.34%  '

$ws.Range('D39').Value = "'0.112"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +13.58%  '

$ws.Range('D40').Value = "'1.35"
$ws.Range('D40').Style = 'Normal'

$ws.Range('D41').Value = "'17.86"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.16%  '

$ws.Range('E43').Value = '  -3.47%  '

$ws.Range('D44').Value = "'96.95"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.87%  '

$ws.Range('E45').Value = '  +1.30%  '

$ws.Range('E46').Value = '  +14.92%  '

$ws.Range('E47').Value = '  +4.69%  '

$ws.Range('D48').Value = "'1.285.51"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.19%  '

$ws.Range('D49').Value = "'2.90"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.78%  '

$ws.Range('D50').Value = "'6.75"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.14%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = "'2.232.58"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.46%  '
